$d = $word.ActiveDocument

# --- Add "red line" (first-line) indentation to the MMain paragraph style ---
# 709 twips ≈ 35.45 pt ≈ 1.25 cm, inserted right after the line-spacing setting.
$mmain = $d.Styles("MMain")
$mmain.ParagraphFormat.FirstLineIndent = 709 / 20

# --- The lone body paragraph (styled MMain) explicitly keeps a zero first-line
#     indent so it is not affected by the new style default. ---
$p = $d.Paragraphs(1)
$p.Range.ParagraphFormat.FirstLineIndent = 0
